$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending clusters (column A) - in order of row blocks
$senders = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")
# Target clusters (column D) - in order within each block
$targets = @("Inflammatory-Mac", "Neutrophils", "Resolving-Mac")

# Numeric data for columns E..T for each of the 12 data rows (row 2 .. row 13)
$data = @(
    @("2", "1", "10.317597", "20.635194", "0.02411218825630575", "0.01628546888406147", "3", "1", "0.3657716666666667", "1.097315", "0.03059585711603819", "0.03059585711603819", "3.773884650685", "22.64330790411", "0.0007377330666449447", "0.0004982678790444307"),
    @("2", "1", "10.317597", "20.635194", "0.02411218825630575", "0.01628546888406147", "3", "1", "11.510206", "34.530618", "0.9627990635838353", "0.9627990635838353", "118.757666894982", "712.546001369892", "0.02321519227412833", "0.01567963419159808"),
    @("2", "1", "10.317597", "20.635194", "0.02411218825630575", "0.01628546888406147", "3", "1", "0.07896333333333333", "0.23689", "0.006605079300126477", "0.006605079300126477", "0.81471185111", "4.88827110666", "0.0001592629155324778", "0.0001075668134189683"),
    @("3", "1", "410.9308523333334", "1232.792557", "0.9603439707700846", "0.9729302679454375", "3", "1", "0.3657716666666667", "1.097315", "0.03059585711603819", "0.03059585711603819", "150.3068627427172", "1352.761764684455", "0.02938254691193026", "0.02976763546192735"),
    @("3", "1", "410.9308523333334", "1232.792557", "0.9603439707700846", "0.9729302679454375", "3", "1", "11.510206", "34.530618", "0.9627990635838353", "0.9627990635838353", "4729.898762112248", "42569.08885901023", "0.9246182757758196", "0.9367363509102371"),
    @("3", "1", "410.9308523333334", "1232.792557", "0.9603439707700846", "0.9729302679454375", "3", "1", "0.07896333333333333", "0.23689", "0.006605079300126477", "0.006605079300126477", "32.44846986974778", "292.03622882773", "0.006343148082334753", "0.006426281573272916"),
    @("2", "1", "6.2889545", "12.577909", "0.01469726476420248", "0.009926591707645529", "3", "1", "0.3657716666666667", "1.097315", "0.03059585711603819", "0.03059585711603819", "2.300321369055833", "13.801928214335", "0.0004496754127221216", "0.0003037125815363721"),
    @("2", "1", "6.2889545", "12.577909", "0.01469726476420248", "0.009926591707645529", "3", "1", "11.510206", "34.530618", "0.9627990635838353", "0.9627990635838353", "72.387161819627", "434.3229709177619", "0.01415051275221784", "0.00955731320070018"),
    @("2", "1", "6.2889545", "12.577909", "0.01469726476420248", "0.009926591707645529", "3", "1", "0.07896333333333333", "0.23689", "0.006605079300126477", "0.006605079300126477", "0.4965968105016667", "2.97958086301", "9.707659926251204E-05", "6.556592540897662E-05"),
    @("2", "0.6666666666666666", "0.3622496666666666", "1.086749", "0.0008465762094071585", "0.0008576714628553975", "3", "1", "0.3657716666666667", "1.097315", "0.03059585711603819", "0.03059585711603819", "0.1325006643261111", "1.192505978935", "2.590172474085865E-05", "2.62411935300272E-05"),
    @("2", "0.6666666666666666", "0.3622496666666666", "1.086749", "0.0008465762094071585", "0.0008576714628553975", "3", "1", "11.510206", "34.530618", "0.9627990635838353", "0.9627990635838353", "4.169568286764666", "37.52611458088199", "0.0008150827816695651", "0.0008257652812997549"),
    @("2", "0.6666666666666666", "0.3622496666666666", "1.086749", "0.0008465762094071585", "0.0008576714628553975", "3", "1", "0.07896333333333333", "0.23689", "0.006605079300126477", "0.006605079300126477", "0.02860444117888888", "0.25743997061", "5.591702996734761E-06", "5.664988025615381E-06")
)

$row = 2
$dataIndex = 0
foreach ($sender in $senders) {
    foreach ($target in $targets) {
        $ws.Cells.Item($row, 1).Value = $sender
        $ws.Cells.Item($row, 2).Value = "Cxcl5"
        $ws.Cells.Item($row, 3).Value = "Cxcr1"
        $ws.Cells.Item($row, 4).Value = $target

        $vals = $data[$dataIndex]
        for ($c = 0; $c -lt $vals.Length; $c++) {
            $ws.Cells.Item($row, 5 + $c).Value = [double]$vals[$c]
        }

        $row++
        $dataIndex++
    }
}
